$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '25.767.98'
$ws.Range("E2").Value = '  +0.30%  '

# Row 3
$ws.Range("D3").Value = '1.746.02'
$ws.Range("E3").Value = '  +0.03%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.19%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.36'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.68%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.07%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5073'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.31%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '40.48'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.59%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2660'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +6.68%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06173'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.24%  '

# Row 11
$ws.Range("D11").Value = '1.751.65'
$ws.Range("E11").Value = '  +0.39%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.06934'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.95%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '15.35'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.80%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6216'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +10.46%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.464'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.17%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '77.41'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.50%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.003'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.22%  '

# Row 18
$ws.Range("E18").Value = '  -0.01%  '

# Row 19
$ws.Range("D19").Value = '25.782.00'
$ws.Range("E19").Value = '  +0.22%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.59'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.85%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.000006634'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.78%  '

# Row 22
$ws.Range("D22").Value = '1.977.22'
$ws.Range("E22").Value = '  +0.50%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.048'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.69%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.241'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.06%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.123'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.20%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '136.73'
$ws.Range("D26").Style = "Normal"

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.455'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.03%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.06'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.13%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.739'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.48%  '

# Row 30
$ws.Range("E30").Value = '  +0.89%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08181'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.13%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.687'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.34%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.390'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.90%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04401'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.37%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.652'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.58%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9904'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.09%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5988'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.10%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.571'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.94%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01560'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.36%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.930'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.42%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.002'
$ws.Range("D41").Style = "Normal"

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '101.37'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.68%  '

# Row 43
$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.3812'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.03%  '

# Row 44
$ws.Range("B44").Value = 'TrustWalletToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.7457'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.40%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.870'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.25%  '

# Row 46
$ws.Range("E46").Value = '  +5.58%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1094'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.29%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.912'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.26%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '29.98'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.27%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '52.51'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.32%  '

# Row 51
$ws.Range("E51").Value = '  +0.50%  '
